# Update "想去人数" (interest count) figures in column F across the
# relevant worksheets, matching the regenerated gh-pages data output.

$wb = $excel.ActiveWorkbook

# Sheet 1: 展览 (Exhibition)
$wsExhibition = $wb.Worksheets.Item(1)
$wsExhibition.Range("F2").Value = 387
$wsExhibition.Range("F3").Value = 1011
$wsExhibition.Range("F4").Value = 245
$wsExhibition.Range("F5").Value = 1371
$wsExhibition.Range("F6").Value = 8444
$wsExhibition.Range("F7").Value = 53
$wsExhibition.Range("F12").Value = 3391
$wsExhibition.Range("F15").Value = 52
$wsExhibition.Range("F16").Value = 924
$wsExhibition.Range("F18").Value = 1090
$wsExhibition.Range("F20").Value = 150
$wsExhibition.Range("F21").Value = 1998

# Sheet 2: 演出 (Performance)
$wsPerformance = $wb.Worksheets.Item(2)
$wsPerformance.Range("F2").Value = 32

# Sheet 4: 全部类型 (All Types)
$wsAll = $wb.Worksheets.Item(4)
$wsAll.Range("F2").Value = 387
$wsAll.Range("F3").Value = 1011
$wsAll.Range("F4").Value = 245
$wsAll.Range("F5").Value = 1371
$wsAll.Range("F6").Value = 8444
$wsAll.Range("F7").Value = 53
$wsAll.Range("F12").Value = 3391
$wsAll.Range("F15").Value = 52
$wsAll.Range("F16").Value = 924
$wsAll.Range("F18").Value = 1090
$wsAll.Range("F20").Value = 150
$wsAll.Range("F21").Value = 1998
$wsAll.Range("F22").Value = 32
